# Update "想去人数" (want-to-go count) column F figures to reflect the
# latest scrape (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 8249
$ws.Range("F4").Value = 2708
$ws.Range("F5").Value = 947
$ws.Range("F6").Value = 305
$ws.Range("F7").Value = 814
$ws.Range("F8").Value = 633
$ws.Range("F9").Value = 110
$ws.Range("F11").Value = 380
$ws.Range("F12").Value = 877
$ws.Range("F13").Value = 3549
$ws.Range("F14").Value = 238
$ws.Range("F15").Value = 130
$ws.Range("F16").Value = 768
$ws.Range("F17").Value = 764
$ws.Range("F19").Value = 475
$ws.Range("F20").Value = 1
$ws.Range("F22").Value = 452
$ws.Range("F23").Value = 1338
$ws.Range("F24").Value = 381
$ws.Range("F26").Value = 139
$ws.Range("F27").Value = 133
$ws.Range("F28").Value = 309
$ws.Range("F29").Value = 41
$ws.Range("F32").Value = 509
$ws.Range("F33").Value = 609
$ws.Range("F34").Value = 35
$ws.Range("F36").Value = 42
$ws.Range("F37").Value = 25
$ws.Range("F38").Value = 227
$ws.Range("F39").Value = 116
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 1
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 8249
$ws.Range("F6").Value = 2708
$ws.Range("F7").Value = 947
$ws.Range("F8").Value = 305
$ws.Range("F9").Value = 814
$ws.Range("F10").Value = 633
$ws.Range("F11").Value = 110
$ws.Range("F13").Value = 380
$ws.Range("F14").Value = 877
$ws.Range("F16").Value = 3550
$ws.Range("F17").Value = 238
$ws.Range("F18").Value = 130
$ws.Range("F20").Value = 768
$ws.Range("F21").Value = 764
$ws.Range("F24").Value = 475
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("F28").Value = 453
$ws.Range("F29").Value = 1338
$ws.Range("F30").Value = 381
$ws.Range("F32").Value = 139
$ws.Range("F33").Value = 133
$ws.Range("F35").Value = 309
$ws.Range("F36").Value = 41
$ws.Range("F39").Value = 509
$ws.Range("F40").Value = 609
$ws.Range("F41").Value = 35
$ws.Range("F43").Value = 42
$ws.Range("F44").Value = 25
$ws.Range("F45").Value = 227
$ws.Range("F46").Value = 116
